$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update K column (최종점수 / final score) for rows 2-5
$ws.Range("K2").Value = 57.3
$ws.Range("K3").Value = 55.5
$ws.Range("K4").Value = 54.3
$ws.Range("K5").Value = 54.3

# Update N column (MACRO_SCORE) for rows 2-5
$ws.Range("N2").Value = 51.15965480231979
$ws.Range("N3").Value = 51.15965480231979
$ws.Range("N4").Value = 51.15965480231979
$ws.Range("N5").Value = 51.15965480231979
